$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Highlight existing homework-score ranges with a green fill (new fill/style entries) ---
$ws.Range("D5:J5").Interior.Color = 5296274
$ws.Range("L5").Interior.Color = 5296274
$ws.Range("D6:G6").Interior.Color = 5296274

# --- New cell K6: same border as the other "д8" score cells, plus the new green fill ---
$ws.Range("L8").Copy()
$ws.Range("K6").PasteSpecial(-4122) | Out-Null
$ws.Range("K6").Value = 5
$ws.Range("K6").Interior.Color = 5296274

# --- Fill in previously empty homework scores (formatting stays as-is) ---
$ws.Range("G11").Value = 5
$ws.Range("H11").Value = 5
$ws.Range("I11").Value = 5

$ws.Range("E14").Value = 5

$ws.Range("I31").Value = 5
$ws.Range("L31").Value = 5

# --- New cell K14: copy formatting from L14 (same style used on that row) ---
$ws.Range("L14").Copy()
$ws.Range("K14").PasteSpecial(-4122) | Out-Null
$ws.Range("K14").Value = 5

# --- New cell J31: copy formatting from J32 (same column/style pattern) ---
$ws.Range("J32").Copy()
$ws.Range("J31").PasteSpecial(-4122) | Out-Null
$ws.Range("J31").Value = 5

# --- New cell K31: copy formatting from K10 (same column/style pattern) ---
$ws.Range("K10").Copy()
$ws.Range("K31").PasteSpecial(-4122) | Out-Null
$ws.Range("K31").Value = 5

# --- New plain value cells in the "Варианты" helper column T ---
$ws.Range("T10").Value = 5
$ws.Range("T19").Value = 5
$ws.Range("T20").Value = 5
$ws.Range("T21").Value = 5

# --- Update the active selection/view to cell K6 ---
$ws.Range("K6").Select() | Out-Null

Write-Host "edits applied"
